$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2-37: 45656 -> 45657
for ($r = 2; $r -le 37; $r++) {
    $ws.Cells.Item($r, 3).Value = 45657
}

# Swap row 36 and row 37 values for columns A (Beteckning) and G (Area (ha))
$ws.Range("A36").Value = "A 60501-2024"
$ws.Range("A37").Value = "A 60500-2024"

$ws.Range("G36").Value = 0.6
$ws.Range("G37").Value = 0.8
